# Commit: "change flow of options download and calucaiton"
#
# The Timesheet_Shashank sheet's trailing rows (Options IV / DB work items)
# are updated: the "Options IV" entry on 2018-07-26 (row 70) is renamed to
# "Options IV and DB", and two previously-blank rows (73 and 74) are filled
# in with new task descriptions and durations. The view's selection is also
# moved from H69 to G65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 70 (2018-07-30): rename the task description.
$ws.Range("C70").Value = "Options IV and DB"

# Row 73 (2018-08-02): fill in task description and duration (hours).
$ws.Range("C73").Value = "DB pushed all stock data, fetching latest date and fetch delta data"
$ws.Range("E73").Value = 6

# Row 74 (2018-08-03): fill in task description and duration (hours).
$ws.Range("C74").Value = "DB, upload delta data,change flow of file Options download and calcultate IV"
$ws.Range("E74").Value = 3

# Update the sheet's active selection to match the saved view state.
[void]$ws.Range("G65").Select()
